# Update "want-to-go" head-count figures (column F) on both the
# "展览" and "全部类型" worksheets, which hold duplicate data tables.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3056
    7  = 1672
    12 = 1375
    16 = 33
    21 = 90
    23 = 3207
    24 = 392
    25 = 137
    26 = 318
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
